# Ajuste para conectarse a DB RDS
#
# 1. "Java" dependency bullet: drop the trailing Red Hat OpenJDK
#    documentation link (and the line break that introduced it), and
#    mark the remaining how2shout Java install link with a leading "*".
# 2. "Git." dependency bullet: add a line break followed by a new
#    how2shout Git install tutorial link.
# 3. "Maven." dependency bullet: add a line break followed by a new
#    LinkedIn Maven setup tutorial link.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Java bullet - remove the Red Hat link, prefix the
# how2shout link with "*".
# ---------------------------------------------------------------------
$javaUrlText = "https://linux.how2shout.com/how-to-install-java-17-on-amazon-linux-2023/"
$redhatUrlText = "https://access.redhat.com/documentation/es-es/red_hat_build_of_openjdk/17/html/installing_and_using_red_hat_build_of_openjdk_17_on_rhel/installing-openjdk11-on-rhel8_openjdk"

# Locate the end of the Java URL and the end of the Red Hat URL so we
# can drop everything in between (the "<w:br/>" plus the Red Hat run).
$javaRng = $d.Content
$javaRng.Find.Execute($javaUrlText, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$javaUrlEnd = $javaRng.End

$redhatRng = $d.Content
$redhatRng.Find.Execute($redhatUrlText, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$redhatUrlEnd = $redhatRng.End

$d.Range($javaUrlEnd, $redhatUrlEnd).Delete()

# Prefix the (now last) Java URL with "*".
$javaRng2 = $d.Content
$javaRng2.Find.Execute($javaUrlText, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$d.Range($javaRng2.Start, $javaRng2.Start).InsertBefore("*")

# ---------------------------------------------------------------------
# Change 2: Git bullet - append a break + the new tutorial link.
# ---------------------------------------------------------------------
$gitRng = $d.Content
$gitRng.Find.Execute("Git.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$gitEnd = $gitRng.End
$d.Range($gitEnd, $gitEnd).InsertAfter( `
    "`vhttps://linux.how2shout.com/how-to-install-git-on-aws-ec2-amazon-linux-2/")

# ---------------------------------------------------------------------
# Change 3: Maven bullet - append a break + the new tutorial link.
# ---------------------------------------------------------------------
$mavenRng = $d.Content
$mavenRng.Find.Execute("Maven.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$mavenEnd = $mavenRng.End
$d.Range($mavenEnd, $mavenEnd).InsertAfter( `
    "`vhttps://www.linkedin.com/pulse/setting-up-maven-aws-ec2-lionel-tchami-nfada-bsc-msc-/")
